$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 833, pushing existing data
# (old rows 833-886) down to 835-888.
$ws.Range("A833:A834").EntireRow.Insert()

# New row 833: Mercado Mayorista Lo Valledor de Santiago, Coliflor, Primera
$ws.Cells.Item(833,1).Value = 6
$ws.Cells.Item(833,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(833,3).Value = "Metropolitana"
$ws.Cells.Item(833,4).Value = 44826
$ws.Cells.Item(833,5).Value = 13
$ws.Cells.Item(833,6).Value = 100112008
$ws.Cells.Item(833,7).Value = "Coliflor"
$ws.Cells.Item(833,8).Value = "Sin especificar"
$ws.Cells.Item(833,9).Value = "Primera"
$ws.Cells.Item(833,10).Value = 8200
$ws.Cells.Item(833,11).Value = 800
$ws.Cells.Item(833,12).Value = 850
$ws.Cells.Item(833,13).Value = 821
$ws.Cells.Item(833,14).Value = "`$/unidad"
$ws.Cells.Item(833,15).Value = "Región Metropolitana"
$ws.Cells.Item(833,16).Value = 821
$ws.Cells.Item(833,17).Value = 1
$ws.Cells.Item(833,18).Value = "Hortaliza"

# New row 834: Mercado Mayorista Lo Valledor de Santiago, Coliflor, Segunda
$ws.Cells.Item(834,1).Value = 6
$ws.Cells.Item(834,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(834,3).Value = "Metropolitana"
$ws.Cells.Item(834,4).Value = 44826
$ws.Cells.Item(834,5).Value = 13
$ws.Cells.Item(834,6).Value = 100112008
$ws.Cells.Item(834,7).Value = "Coliflor"
$ws.Cells.Item(834,8).Value = "Sin especificar"
$ws.Cells.Item(834,9).Value = "Segunda"
$ws.Cells.Item(834,10).Value = 3100
$ws.Cells.Item(834,11).Value = 700
$ws.Cells.Item(834,12).Value = 700
$ws.Cells.Item(834,13).Value = 700
$ws.Cells.Item(834,14).Value = "`$/unidad"
$ws.Cells.Item(834,15).Value = "Región Metropolitana"
$ws.Cells.Item(834,16).Value = 700
$ws.Cells.Item(834,17).Value = 1
$ws.Cells.Item(834,18).Value = "Hortaliza"
